{"js": "// Apply the dated worksheet refresh: update the header date and every\n// three-digit-by-one-digit multiplication prompt in the table.\nconst replacements = [\n  [\"2024-08-27 Tuesday\", \"2024-08-28 Wednesday\"],\n  [\"374\u00d73=\", \"661\u00d78=\"],\n  [\"104\u00d73=\", \"737\u00d74=\"],\n  [\"198\u00d73=\", \"469\u00d79=\"],\n  [\"251\u00d75=\", \"790\u00d77=\"],\n  [\"463\u00d75=\", \"720\u00d78=\"],\n  [\"550\u00d75=\", \"763\u00d73=\"],\n  [\"387\u00d75=\", \"964\u00d74=\"],\n  [\"554\u00d76=\", \"693\u00d78=\"],\n  [\"604\u00d77=\", \"202\u00d72=\"],\n  [\"783\u00d77=\", \"552\u00d76=\"],\n  [\"974\u00d75=\", \"590\u00d77=\"],\n  [\"110\u00d72=\", \"844\u00d73=\"],\n  [\"311\u00d74=\", \"551\u00d74=\"],\n  [\"115\u00d74=\", \"950\u00d77=\"],\n  [\"542\u00d76=\", \"145\u00d79=\"],\n  [\"964\u00d79=\", \"771\u00d73=\"],\n  [\"414\u00d73=\", \"217\u00d75=\"],\n  [\"230\u00d77=\", \"357\u00d75=\"],\n  [\"826\u00d75=\", \"345\u00d75=\"],\n  [\"120\u00d75=\", \"226\u00d73=\"],\n  [\"419\u00d78=\", \"122\u00d75=\"],\n  [\"132\u00d73=\", \"943\u00d73=\"],\n  [\"554\u00d74=\", \"139\u00d79=\"],\n  [\"837\u00d74=\", \"138\u00d73=\"],\n  [\"555\u00d73=\", \"584\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the dated worksheet refresh: update the header date and every\n# three-digit-by-one-digit multiplication prompt in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-27 Tuesday\", \"2024-08-28 Wednesday\"),\n    @(\"374\u00d73=\", \"661\u00d78=\"),\n    @(\"104\u00d73=\", \"737\u00d74=\"),\n    @(\"198\u00d73=\", \"469\u00d79=\"),\n    @(\"251\u00d75=\", \"790\u00d77=\"),\n    @(\"463\u00d75=\", \"720\u00d78=\"),\n    @(\"550\u00d75=\", \"763\u00d73=\"),\n    @(\"387\u00d75=\", \"964\u00d74=\"),\n    @(\"554\u00d76=\", \"693\u00d78=\"),\n    @(\"604\u00d77=\", \"202\u00d72=\"),\n    @(\"783\u00d77=\", \"552\u00d76=\"),\n    @(\"974\u00d75=\", \"590\u00d77=\"),\n    @(\"110\u00d72=\", \"844\u00d73=\"),\n    @(\"311\u00d74=\", \"551\u00d74=\"),\n    @(\"115\u00d74=\", \"950\u00d77=\"),\n    @(\"542\u00d76=\", \"145\u00d79=\"),\n    @(\"964\u00d79=\", \"771\u00d73=\"),\n    @(\"414\u00d73=\", \"217\u00d75=\"),\n    @(\"230\u00d77=\", \"357\u00d75=\"),\n    @(\"826\u00d75=\", \"345\u00d75=\"),\n    @(\"120\u00d75=\", \"226\u00d73=\"),\n    @(\"419\u00d78=\", \"122\u00d75=\"),\n    @(\"132\u00d73=\", \"943\u00d73=\"),\n    @(\"554\u00d74=\", \"139\u00d79=\"),\n    @(\"837\u00d74=\", \"138\u00d73=\"),\n    @(\"555\u00d73=\", \"584\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
